$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A53").Value = 52
$ws.Range("B53").Value = 84
$ws.Range("C53").Value = 1
$ws.Range("D53").Value = 11
$ws.Range("E53").Value = 22
$ws.Range("F53").Value = 96
$ws.Range("G53").Value = 118
